# Update the SFIA skills table (Sheet1, A1:D34) to reflect the new job spec
# (SWDN, SINT, TEST, CFMG, REQM skills replacing/augmenting PROG/MADE/CSOP rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@("SkillCode", "SFIA Level", "Keycode", "Description")
    ,@("Autonomy", 1, 1, "Works under close direction")
    ,@("Autonomy", 1, 2, "Uses little discretion in attending to enquiries")
    ,@("Autonomy", 1, 3, "Is expected to seek guidance in unexpected situations")
    ,@("Influence", 1, 1, "Minimal influence")
    ,@("Influence", 1, 2, "May work alone or interact with immediate colleagues")
    ,@("Complexity", 1, 1, "Performs routine activities in a structured environment")
    ,@("Complexity", 1, 2, "Requires assistance in resolving unexpected problems")
    ,@("Complexity", 1, 3, "Participates in the generation of new ideas")
    ,@("Knowledge", 1, 1, "Has a basic generic knowledge appropriate to area of work")
    ,@("Knowledge", 1, 2, "Applies newly acquired knowledge to develop new skills")
    ,@("PROG", 1, 1, "Designs, codes, verifies, tests, documents, amends and refactors simple programs/scripts")
    ,@("PROG", 1, 2, "Applies agreed standards and tools to achieve a well-engineered result")
    ,@("PROG", 1, 3, "Reviews own work")
    ,@("SWDN", 1, 4, "Reviews own work")
    ,@("PROG", 1, 1, "Designs, codes, verifies, tests, documents, amends and refactors simple programs/scripts")
    ,@("PROG", 1, 2, "Applies agreed standards and tools to achieve a well-engineered result")
    ,@("PROG", 1, 3, "Reviews own work")
    ,@("SINT", 1, 1, "Produces software builds from software source code")
    ,@("SINT", 1, 2, "Conducts tests as defined in an integration test specification and records the details of any failures")
    ,@("SINT", 1, 3, "Analyses and reports on integration test activities and results")
    ,@("SINT", 1, 4, "Identifies and reports issues and risks")
    ,@("TEST", 1, 1, "Designs test cases, creates test scripts and test data, and automates repeatable tasks working to the requirements or specifications provided")
    ,@("TEST", 1, 2, "Defines test conditions for given requirements")
    ,@("TEST", 1, 3, "Executes and records manual and automated testing  in accordance with test plans")
    ,@("TEST", 1, 4, "Analyses and reports on test activities, results, issues and risks")
    ,@("CFMG", 1, 1, "Applies tools, techniques and processes to administer, track, log, report on and correct configuration items, components and changes")
    ,@("CFMG", 1, 2, "Assists with audits to check the accuracy of the information and undertakes any necessary corrective action under direction")
    ,@("MADE", 1, 1, "Example MT")
    ,@($null, $null, $null, $null)
    ,@("REQM", 1, 1, "Uses standard techniques to elicit, specify, and document requirements for simple subject areas with clearly-defined boundaries")
    ,@("REQM", 1, 2, "Assists in the definition and management of requirements")
    ,@("REQM", 1, 3, "Assists in the creation of a requirements baseline")
    ,@("REQM", 1, 4, "Assists in investigating and applying authorised changes")
)

$data = New-Object 'object[,]' $rows.Count,4
for ($r = 0; $r -lt $rows.Count; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        $data[$r, $c] = $rows[$r][$c]
    }
}

$ws.Range("A1:D34").Value = $data
